$p = $ppt.ActivePresentation

# 1) Resize the "Subtitle 2" placeholder on slide layout 1 (Title Slide):
#    cx 11676184 EMU -> 11119442 EMU  (width in points = EMU / 12700)
$master = $p.SlideMaster
$layout1 = $master.CustomLayouts.Item(1)
$subtitle = $layout1.Shapes.Item(2)
$subtitle.Width = 11119442 / 12700

# 2) Update the cached date text of the Date Placeholder on the slide master
$datePh = $master.Shapes.Item(3)
$datePh.TextFrame.TextRange.Text = "10/23/2018"
